# Apply cell updates for the "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced into a number by Excel
# (e.g. "0.110" -> 0.11, "1.00" -> 1) -- force them to Text format first so the
# literal string is preserved exactly as in the source data.
$numericLookingCells = @(
    "D5",
    "D6",
    "D7",
    "D8",
    "D10",
    "D11",
    "D12",
    "D15",
    "D16",
    "D19",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D29",
    "D30",
    "D31",
    "D33",
    "D34",
    "D35",
    "D36",
    "D39",
    "D40",
    "D41",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51",
)
foreach ($ref in $numericLookingCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# --- Cell value assignments (only the cells that actually changed) ---

$ws.Range("D2").Value = '65.303.04'
$ws.Range("E2").Value = '  +3.19%  '
$ws.Range("D3").Value = '3.489.91'
$ws.Range("E3").Value = '  +2.63%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '580.76'
$ws.Range("E5").Value = '  +2.40%  '
$ws.Range("D6").Value = '163.14'
$ws.Range("E6").Value = '  +4.83%  '
$ws.Range("D7").Value = '0.617'
$ws.Range("E7").Value = '  +12.98%  '
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '3.492.78'
$ws.Range("E9").Value = '  +2.77%  '
$ws.Range("D10").Value = '7.25'
$ws.Range("E10").Value = '  -1.89%  '
$ws.Range("D11").Value = '0.126'
$ws.Range("E11").Value = '  +3.37%  '
$ws.Range("D12").Value = '0.449'
$ws.Range("E12").Value = '  +4.01%  '
$ws.Range("D13").Value = '4.085.94'
$ws.Range("E13").Value = '  +2.52%  '
$ws.Range("E14").Value = '  +0.49%  '
$ws.Range("D15").Value = '0.0000194'
$ws.Range("E15").Value = '  +0.70%  '
$ws.Range("D16").Value = '28.78'
$ws.Range("E16").Value = '  +5.96%  '
$ws.Range("D17").Value = '65.269.46'
$ws.Range("E17").Value = '  +3.03%  '
$ws.Range("D18").Value = '3.482.09'
$ws.Range("E18").Value = '  +0.00%  '
$ws.Range("D19").Value = '6.48'
$ws.Range("E19").Value = '  +3.86%  '
$ws.Range("D20").Value = '14.42'
$ws.Range("E20").Value = '  +2.64%  '
$ws.Range("D21").Value = '383.39'
$ws.Range("E21").Value = '  +1.26%  '
$ws.Range("D22").Value = '8.24'
$ws.Range("E22").Value = '  +2.37%  '
$ws.Range("D23").Value = '0.556'
$ws.Range("E23").Value = '  +5.26%  '
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  +0.15%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = '72.72'
$ws.Range("E25").Value = '  +1.82%  '
$ws.Range("D26").Value = '0.0000121'
$ws.Range("E26").Value = '  +0.94%  '
$ws.Range("D27").Value = '10.04'
$ws.Range("E27").Value = '  +6.43%  '
$ws.Range("E28").Value = '  +0.31%  '
$ws.Range("B29").Value = 'Fetch.AI'
$ws.Range("C29").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D29").Value = '1.55'
$ws.Range("E29").Value = '  +13.62%  '
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("D31").Value = '6.19'
$ws.Range("E31").Value = '  +2.59%  '
$ws.Range("E32").Value = '  +3.13%  '
$ws.Range("D33").Value = '23.74'
$ws.Range("E33").Value = '  +2.49%  '
$ws.Range("D34").Value = '7.21'
$ws.Range("E34").Value = '  +6.43%  '
$ws.Range("D35").Value = '1.64'
$ws.Range("E35").Value = '  +14.12%  '
$ws.Range("D36").Value = '162.35'
$ws.Range("E36").Value = '  +1.69%  '
$ws.Range("E37").Value = '  +6.22%  '
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '3.020.74'
$ws.Range("E38").Value = '  +2.12%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '0.0782'
$ws.Range("E39").Value = '  +3.84%  '
$ws.Range("D40").Value = '6.84'
$ws.Range("E40").Value = '  +6.88%  '
$ws.Range("D41").Value = '26.98'
$ws.Range("E41").Value = '  +0.28%  '
$ws.Range("E42").Value = '  +6.07%  '
$ws.Range("D43").Value = '0.0323'
$ws.Range("E43").Value = '  +2.41%  '
$ws.Range("D44").Value = '42.99'
$ws.Range("E44").Value = '  +3.02%  '
$ws.Range("D45").Value = '0.785'
$ws.Range("E45").Value = '  +3.07%  '
$ws.Range("D46").Value = '26.04'
$ws.Range("E46").Value = '  +12.05%  '
$ws.Range("D47").Value = '1.12'
$ws.Range("E47").Value = '  +4.29%  '
$ws.Range("D48").Value = '319.54'
$ws.Range("E48").Value = '  +9.34%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").Value = '0.110'
$ws.Range("E49").Value = '  +6.92%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = '6.73'
$ws.Range("E50").Value = '  +6.19%  '
$ws.Range("D51").Value = '0.877'
$ws.Range("E51").Value = '  +4.99%  '
